# The sheet lists weekly Alcachofa (artichoke) price observations for
# "Macroferia Regional de Talca". A new weekly record was added at the
# top of the Madrigal/Primera block (row 116), pushing the existing
# rows 116-132 down to 117-133.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 116; formatting (incl. the date number format
# on column D) is inherited from the row that used to be at 116.
$ws.Rows.Item(116).Insert()

# Populate the new row with the latest weekly observation.
$ws.Cells.Item(116, 1).Value = 5
$ws.Cells.Item(116, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(116, 3).Value = "Maule"
$ws.Cells.Item(116, 4).Value = 45154
$ws.Cells.Item(116, 5).Value = 7
$ws.Cells.Item(116, 6).Value = 100112013
$ws.Cells.Item(116, 7).Value = "Alcachofa"
$ws.Cells.Item(116, 8).Value = "Madrigal"
$ws.Cells.Item(116, 9).Value = "Primera"
$ws.Cells.Item(116, 10).Value = 300
$ws.Cells.Item(116, 11).Value = 13000
$ws.Cells.Item(116, 12).Value = 13000
$ws.Cells.Item(116, 13).Value = 13000
$ws.Cells.Item(116, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(116, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(116, 16).Value = 325
$ws.Cells.Item(116, 17).Value = 40
$ws.Cells.Item(116, 18).Value = "Hortaliza"
